$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ordered (cell, new value) pairs mirroring the 23-Jan-2024 cryptos refresh.
$updates = @(
    @('D2', '39.057.10'),
    @('E2', '  -4.31%  '),
    @('D3', '2.224.76'),
    @('E3', '  -6.99%  '),
    @('E4', '  -0.02%  '),
    @('D5', '296.86'),
    @('E5', '  -5.18%  '),
    @('D6', '80.67'),
    @('E6', '  -8.14%  '),
    @('D7', '0.508'),
    @('E7', '  -4.34%  '),
    @('E8', '  +0.02%  '),
    @('D9', '0.457'),
    @('E9', '  -7.08%  '),
    @('D10', '0.0774'),
    @('E10', '  -6.40%  '),
    @('D11', '27.97'),
    @('E11', '  -10.08%  '),
    @('D12', '46.58'),
    @('E12', '  -12.44%  '),
    @('D13', '0.107'),
    @('E13', '  -1.28%  '),
    @('D14', '2.561.83'),
    @('E14', '  -7.38%  '),
    @('D15', '6.08'),
    @('E15', '  -9.01%  '),
    @('D16', '13.99'),
    @('E16', '  -7.55%  '),
    @('D17', '2.243.78'),
    @('E17', '  -7.02%  '),
    @('D18', '0.712'),
    @('E18', '  -6.47%  '),
    @('D19', '38.936.16'),
    @('E19', '  -4.49%  '),
    @('D20', '0.0₃0856'),
    @('E20', '  -6.33%  '),
    @('D21', '5.74'),
    @('E21', '  -6.82%  '),
    @('D22', '65.26'),
    @('E22', '  -6.80%  '),
    @('D23', '9.84'),
    @('E23', '  -8.70%  '),
    @('D24', '224.74'),
    @('E24', '  -5.49%  '),
    @('E25', '  +0.10%  '),
    @('D26', '2.38'),
    @('E26', '  -9.78%  '),
    @('D27', '1.73'),
    @('E27', '  -5.85%  '),
    @('D28', '22.22'),
    @('E28', '  -6.55%  '),
    @('E29', '  -1.75%  '),
    @('D30', '8.90'),
    @('E30', '  -5.54%  '),
    @('D31', '147.96'),
    @('E31', '  -5.65%  '),
    @('D32', '31.51'),
    @('E32', '  -7.34%  '),
    @('D33', '1.00'),
    @('E33', '  +0.00%  '),
    @('D34', '4.75'),
    @('E34', '  -9.50%  '),
    @('E35', '  -4.43%  '),
    @('D36', '0.0682'),
    @('E36', '  -6.94%  '),
    @('E37', '  -3.98%  '),
    @('D38', '0.0970'),
    @('E38', '  -2.21%  '),
    @('D39', '2.63'),
    @('E39', '  -6.68%  '),
    @('D40', '14.71'),
    @('E40', '  -7.79%  '),
    @('D41', '1.60'),
    @('E41', '  -8.35%  '),
    @('D42', '3.62'),
    @('E42', '  -5.71%  '),
    @('D43', '1.903.38'),
    @('E43', '  -3.32%  '),
    @('D44', '2.16'),
    @('E44', '  -5.03%  '),
    @('D45', '0.0252'),
    @('E45', '  -6.84%  '),
    @('D46', '16.35'),
    @('E46', '  -8.66%  '),
    @('D47', '8.93'),
    @('E47', '  -3.89%  '),
    @('E48', '  -10.59%  '),
    @('D49', '2.450.34'),
    @('E49', '  -7.15%  '),
    @('D50', '87.28'),
    @('E50', '  -6.78%  '),
    @('D51', '66.36'),
    @('E51', '  -9.53%  ')
)

foreach ($pair in $updates) {
    $cell = $pair[0]
    $value = $pair[1]
    $range = $ws.Range($cell)
    # Numeric-looking strings (e.g. "296.86") would otherwise be
    # auto-converted to a Number by Excel on assignment; force Text
    # so the cell keeps the exact literal the site renders.
    if ($value -match '^[+-]?\d+(\.\d+)?$') {
        $range.NumberFormat = '@'
    }
    $range.Value = $value
}
